$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update project descriptions for Cave Cart (D2) and Offline RL Using Transformers (D8)
# replacing the single space between sentences with a literal "\n" marker,
# matching how the site's CSV-driven description text marks paragraph breaks.
$ws.Range("D2").Value = 'In 2019 I attended a Danish folk high school, where I took the game design course. We were a group of three people on the course. I made several games with my group. Cave Cart was the final project, which took a month to make. My main responsibility on the project was programming using Unity and C#.\nThe game is similar to temple run. The player is driving a mine-cart in a tunnel with three tracks. To avoid obstacles like holes and explosive barrels the player switches tracks. The player gets points according to how long they survive and gets extra points by collecting pickups. The maze is procedurally generated, and the game is endless. The game also contains power-ups, a start screen, setting and credit scene.'
$ws.Range("D8").Value = 'During our fourth semester me and my group at DTU made a large machine learning project. All students on fourth semester chose a project from a list. We chose a project about offline reinforcement learning using transformers. My main focus on the project was programming, but also writing the report. The goal of the project was to replicate the results of a research paper, and build on the research paper to obtain new results.\nOffline RL is the task of training a RL model on existing data without being able to collect new data from the environment. This is generally more difficult than online RL because the model cannot sample new data in places where it lacks knowledge. The researchers proposed that a transformer model could be used for offline RL. This is based on the intuition that RL is about modelling a sequence of actions, and transformers are good at modelling sequences.\nIn the paper the problem was to teach an AI how to walk. That AI could be a simulated spider or humanoid. We used MoJuCo to simulate the environment.\nIn the end the transformer AI was able to walk for about 20 seconds before falling, but didn''t seem to generalize much from the data. That means the model performed well when the environment was similar to the data, but not so well when having to exstrapolate from the data.\nWe exstended on the research by'

# These two cells previously had an explicit wrap-text style (with a custom
# row height to match). Restore them to the default "Normal" style and let
# the row heights return to the sheet's standard height.
$ws.Range("D2").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(8).AutoFit()

# Update the active selection to reflect where the author left the cursor.
$ws.Range("E12").Select() | Out-Null
